$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 58,4
$data[0,0] = 0
$data[0,1] = 1.544980316747408
$data[0,2] = 6.86676254930352
$data[0,3] = 0.3644367224954216
$data[1,0] = 1
$data[1,1] = 3.831280865411461
$data[1,2] = 1.950717946524155
$data[1,3] = 0.3162575504973723
$data[2,0] = 2
$data[2,1] = 6.429665394968044
$data[2,2] = 3.504543068054611
$data[2,3] = 0.3268591321286358
$data[3,0] = 3
$data[3,1] = 7.269686110168726
$data[3,2] = 10.27585520416362
$data[3,3] = 0.4765889838384171
$data[4,0] = 4
$data[4,1] = 8.201026164251651
$data[4,2] = 4.479335480805395
$data[4,3] = 0.07270861129031886
$data[5,0] = 5
$data[5,1] = 10.57499949525755
$data[5,2] = 10.0283787213407
$data[5,3] = 0.2254333055046017
$data[6,0] = 6
$data[6,1] = 10.83456178137682
$data[6,2] = 10.3657065302169
$data[6,3] = 0.3261720398712445
$data[7,0] = 7
$data[7,1] = 11.23094531496357
$data[7,2] = 8.425281253103092
$data[7,3] = 0.4450104157212572
$data[8,0] = 8
$data[8,1] = 12.97413078332409
$data[8,2] = 9.154372486332297
$data[8,3] = 0.5037440713291443
$data[9,0] = 9
$data[9,1] = 13.54664878204631
$data[9,2] = 7.244378549995161
$data[9,3] = 0.3965441495492412
$data[10,0] = 10
$data[10,1] = 14.39266784590064
$data[10,2] = 9.425295021221055
$data[10,3] = 0.3136544282744531
$data[11,0] = 11
$data[11,1] = 14.94142834267681
$data[11,2] = 12.14929210939959
$data[11,3] = 0.259107157619748
$data[12,0] = 12
$data[12,1] = 15.80020168688101
$data[12,2] = 7.454097089100859
$data[12,3] = 0.3596141812573913
$data[13,0] = 13
$data[13,1] = 16.00851040643039
$data[13,2] = 5.070337929725086
$data[13,3] = 0.3385502891614018
$data[14,0] = 14
$data[14,1] = 23.68807458007721
$data[14,2] = 2.629972134090708
$data[14,3] = 0.4419241717517934
$data[15,0] = 15
$data[15,1] = 27.72416160122677
$data[15,2] = 4.574855980161045
$data[15,3] = 0.1766204632143659
$data[16,0] = 16
$data[16,1] = 29.46084561329117
$data[16,2] = 6.34750718472964
$data[16,3] = 0.2928287452789782
$data[17,0] = 17
$data[17,1] = 31.60208960550924
$data[17,2] = 2.303324223296186
$data[17,3] = 0.2736507320835925
$data[18,0] = 18
$data[18,1] = 32.99315184678982
$data[18,2] = 5.74426029281828
$data[18,3] = 0.4391810489871426
$data[19,0] = 19
$data[19,1] = 38.06580449424434
$data[19,2] = 7.469818072101941
$data[19,3] = 0.4006112338370554
$data[20,0] = 20
$data[20,1] = 40.3399765718673
$data[20,2] = 10.80177640397956
$data[20,3] = 0.4011462235920872
$data[21,0] = 21
$data[21,1] = 40.37373718621435
$data[21,2] = 4.588168823523216
$data[21,3] = 0.2533197120128846
$data[22,0] = 22
$data[22,1] = 41.04059076197027
$data[22,2] = 5.802603517274036
$data[22,3] = 0.484812522708932
$data[23,0] = 23
$data[23,1] = 42.13920506587948
$data[23,2] = 6.297712411182865
$data[23,3] = 0.2333653029251783
$data[24,0] = 24
$data[24,1] = 42.16911071793339
$data[24,2] = 7.57016633389415
$data[24,3] = 0.4343481839159654
$data[25,0] = 25
$data[25,1] = 43.51451832428538
$data[25,2] = 1.391823011936882
$data[25,3] = 0.1649099096245734
$data[26,0] = 26
$data[26,1] = 43.67201355429995
$data[26,2] = 5.575678929745554
$data[26,3] = 0.3167272162568531
$data[27,0] = 27
$data[27,1] = 47.23160760964704
$data[27,2] = 5.037940506386185
$data[27,3] = 0.1747182997935385
$data[28,0] = 28
$data[28,1] = 49.76220711470757
$data[28,2] = 4.681979678798008
$data[28,3] = 0.3685171341532826
$data[29,0] = 29
$data[29,1] = 52.81700867667762
$data[29,2] = 8.25339910049318
$data[29,3] = 0.3156111377998286
$data[30,0] = 30
$data[30,1] = 57.38481985170372
$data[30,2] = 6.226892358625991
$data[30,3] = 0.1602293672015526
$data[31,0] = 31
$data[31,1] = 62.47454224640119
$data[31,2] = 4.017554352513008
$data[31,3] = 0.2475842874905581
$data[32,0] = 32
$data[32,1] = 63.1359781402728
$data[32,2] = 3.847240142695312
$data[32,3] = 0.3641734560947129
$data[33,0] = 33
$data[33,1] = 63.92922427265702
$data[33,2] = 4.893296615131952
$data[33,3] = 0.1512289517652328
$data[34,0] = 34
$data[34,1] = 63.94030168035091
$data[34,2] = 5.421987803397015
$data[34,3] = 0.15485434331586
$data[35,0] = 35
$data[35,1] = 68.14656561211218
$data[35,2] = 13.22695078881117
$data[35,3] = 0.2904885786904932
$data[36,0] = 36
$data[36,1] = 70.94486736385976
$data[36,2] = 4.013735753902317
$data[36,3] = 0.1689296847075247
$data[37,0] = 37
$data[37,1] = 71.19194662457122
$data[37,2] = 10.73960845972079
$data[37,3] = 0.3340334375894475
$data[38,0] = 38
$data[38,1] = 71.48295976691124
$data[38,2] = 6.207717097863986
$data[38,3] = 0.5114993753167597
$data[39,0] = 39
$data[39,1] = 73.64355284419571
$data[39,2] = 7.202163986162057
$data[39,3] = 0.1770488107516146
$data[40,0] = 40
$data[40,1] = 75.13951730112467
$data[40,2] = 4.995257262382911
$data[40,3] = 0.3641078006113884
$data[41,0] = 41
$data[41,1] = 75.22285217669348
$data[41,2] = 6.825453730610984
$data[41,3] = 0.4539342361253993
$data[42,0] = 42
$data[42,1] = 75.7327716880055
$data[42,2] = 3.483156946358522
$data[42,3] = 0.2605950474001671
$data[43,0] = 43
$data[43,1] = 76.4937900612946
$data[43,2] = 2.248643878161143
$data[43,3] = 0.1684841110986813
$data[44,0] = 44
$data[44,1] = 76.94040150519551
$data[44,2] = 5.681178255061339
$data[44,3] = 0.1924233703323319
$data[45,0] = 45
$data[45,1] = 77.73959333018695
$data[45,2] = 6.566928344840817
$data[45,3] = 0.2882438158694063
$data[46,0] = 46
$data[46,1] = 78.55187535298148
$data[46,2] = 4.419517394600299
$data[46,3] = 0.1009281509786034
$data[47,0] = 47
$data[47,1] = 79.23302284632386
$data[47,2] = 7.779453545033725
$data[47,3] = 0.4715775145352027
$data[48,0] = 48
$data[48,1] = 81.55021290463674
$data[48,2] = 13.66891840741803
$data[48,3] = 0.2763729393639068
$data[49,0] = 49
$data[49,1] = 85.74064320056776
$data[49,2] = 6.970726387345792
$data[49,3] = 0.4509438091587586
$data[50,0] = 50
$data[50,1] = 89.14545095999485
$data[50,2] = 9.849588575994993
$data[50,3] = 0.465543040083724
$data[51,0] = 51
$data[51,1] = 89.58494029668941
$data[51,2] = 7.375291471533125
$data[51,3] = 0.3848981083860228
$data[52,0] = 52
$data[52,1] = 94.7063218580974
$data[52,2] = 5.971605045621442
$data[52,3] = 0.3456471090593382
$data[53,0] = 53
$data[53,1] = 96.05248213573093
$data[53,2] = 7.201848912652274
$data[53,3] = 0.5795711998767099
$data[54,0] = 54
$data[54,1] = 98.04086959591393
$data[54,2] = 4.687995252470976
$data[54,3] = 0.4387299171901498
$data[55,0] = 55
$data[55,1] = 98.42020142694751
$data[55,2] = 11.89499501616148
$data[55,3] = 0.4242084671186639
$data[56,0] = 56
$data[56,1] = 99.44767940872428
$data[56,2] = 5.701716566501312
$data[56,3] = 0.5279201860832713
$data[57,0] = 57
$data[57,1] = 99.97122813223787
$data[57,2] = 10.62120113248827
$data[57,3] = 0.4286642951985921

# Write existing + new data rows (2 through 59) in one pass
$ws.Range("A2:D59").Value = $data

# New rows (54-59) need the same formatting (bold, border, centered/top aligned)
# that the other A-column index cells use. Copy format from A2 and paste it
# onto the newly added index cells so the style matches exactly.
$ws.Range("A2").Copy()
$ws.Range("A54:A59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
